# Refresh cryptos list data: updated prices / 1h volume percentages,
# and swap the Hedera/Kaspa row order (rows 31-32), per the upstream
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'62.161.85"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "'3.435.88"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'409.01"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'129.45"
$ws.Range("E6").Value = "  -4.49%  "
$ws.Range("D7").Value = "'0.634"
$ws.Range("E7").Value = "  +6.84%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.757"
$ws.Range("E9").Value = "  +11.04%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("E10").Value = "  +15.98%  "
$ws.Range("D11").Value = "'43.00"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'0.140"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").Value = "'20.39"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "'0.0000191"
$ws.Range("E15").Value = "  +48.35%  "
$ws.Range("D16").Value = "'3.438.97"
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "'62.180.63"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "'11.36"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").Value = "'389.80"
$ws.Range("E20").Value = "  +24.25%  "
$ws.Range("D21").Value = "'88.41"
$ws.Range("E21").Value = "  +4.78%  "
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").Value = "'13.34"
$ws.Range("E23").Value = "  +4.71%  "
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("D25").Value = "'32.04"
$ws.Range("E25").Value = "  +8.45%  "
$ws.Range("D26").Value = "'4.80"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "'8.47"
$ws.Range("E27").Value = "  +1.30%  "
$ws.Range("D28").Value = "'7.68"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'2.76"
$ws.Range("E29").Value = "  +11.47%  "
$ws.Range("D30").Value = "'44.08"
$ws.Range("E30").Value = "  +6.51%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "'0.171"
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("D33").Value = "'11.79"
$ws.Range("E33").Value = "  +4.01%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.0492"
$ws.Range("E35").Value = "  +1.85%  "
$ws.Range("D36").Value = "'52.26"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  -2.53%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +6.73%  "
$ws.Range("D41").Value = "'0.312"
$ws.Range("E41").Value = "  +8.15%  "
$ws.Range("D42").Value = "'141.65"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "'1.96"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'3.99"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").Value = "'16.69"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").Value = "'21.86"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("D48").Value = "'2.115.89"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").Value = "'1.93"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  +6.08%  "
